$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Danh sach suat chie^'u" - merge the split runs "chie^" + "u" into
#    a single run "chieu" (paragraph: "Output: Danh sach suat chieu").
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(18)
$p1.Range.Find.Execute("chiếu", $false, $false, $false, $false, $false, `
    $true, 1, $false, "chiếu", 2)

# ---------------------------------------------------------------------
# 2) "Kiem tra (ID_PHONG, ID_PHIM, ID_CA_CHIEU) duy nhat" - merge the
#    split runs " " + "(ID_PHONG, ID_PHIM, ID_CA_CHIEU) " into one run.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(28)
$p2.Range.Find.Execute(" (ID_PHONG, ID_PHIM, ID_CA_CHIEU) ", $false, $false, `
    $false, $false, $false, $true, 1, $false, " (ID_PHONG, ID_PHIM, ID_CA_CHIEU) ", 2)

# ---------------------------------------------------------------------
# 3) "Phong" bullet under "Xoa" - merge split runs "P" + "hong" into
#    a single run "Phong".
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(56)
$p3.Range.Find.Execute("Phòng", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Phòng", 2)

# ---------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the end of the first "Sua (x,y)"
#    paragraph (under So_Ghe) into the middle of "ID_LOAI_PHONG)" in the
#    "Sua (ID_LOAI_PHONG)" paragraph, splitting that run into
#    "ID_LO" + bookmark + "AI_PHONG)".
# ---------------------------------------------------------------------
$pTarget = $d.Paragraphs(58)
$rFind = $pTarget.Range.Duplicate
$rFind.Find.Execute("ID_LO")
$insertPoint = $rFind.End

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$bmRange = $d.Range($insertPoint, $insertPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 5) Apply strikethrough formatting to the second "Sua (x,y)" paragraph
#    (under So_day), which previously had no strike formatting.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(68)
$p5.Range.Font.StrikeThrough = $true
